$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric columns (E:T) for existing rows 2-13 with refreshed TPM-based values
# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.767552
$ws.Range("H2").Value = 8.302655999999999
$ws.Range("I2").Value = 0.04706493447833917
$ws.Range("J2").Value = 0.04706493447833917
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7489546666666667
$ws.Range("N2").Value = 2.246864
$ws.Range("O2").Value = 0.05220789806691288
$ws.Range("P2").Value = 0.05220789806691287
$ws.Range("Q2").Value = 2.072770985642666
$ws.Range("R2").Value = 18.654938870784
$ws.Range("S2").Value = 0.002457161301771065
$ws.Range("T2").Value = 0.002457161301771064

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.767552
$ws.Range("H3").Value = 8.302655999999999
$ws.Range("I3").Value = 0.04706493447833917
$ws.Range("J3").Value = 0.04706493447833917
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.843693666666667
$ws.Range("N3").Value = 8.531081
$ws.Range("O3").Value = 0.1982273102638064
$ws.Range("P3").Value = 0.1982273102638064
$ws.Range("Q3").Value = 7.870070094570666
$ws.Range("R3").Value = 70.830630851136
$ws.Range("S3").Value = 0.00932955536938346
$ws.Range("T3").Value = 0.009329555369383458

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.767552
$ws.Range("H4").Value = 8.302655999999999
$ws.Range("I4").Value = 0.04706493447833917
$ws.Range("J4").Value = 0.04706493447833917
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.337765
$ws.Range("N4").Value = 31.013295
$ws.Range("O4").Value = 0.7206216949842531
$ws.Range("P4").Value = 0.720621694984253
$ws.Range("Q4").Value = 28.61030220128
$ws.Range("R4").Value = 257.49271981152
$ws.Range("S4").Value = 0.03391601285810358
$ws.Range("T4").Value = 0.03391601285810358

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.767552
$ws.Range("H5").Value = 8.302655999999999
$ws.Range("I5").Value = 0.04706493447833917
$ws.Range("J5").Value = 0.04706493447833917
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4152066666666667
$ws.Range("N5").Value = 1.24562
$ws.Range("O5").Value = 0.02894309668502767
$ws.Range("P5").Value = 0.02894309668502767
$ws.Range("Q5").Value = 1.149106040746667
$ws.Range("R5").Value = 10.34195436672
$ws.Range("S5").Value = 0.001362204949081063
$ws.Range("T5").Value = 0.001362204949081063

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 21.05317333333333
$ws.Range("H6").Value = 63.15952
$ws.Range("I6").Value = 0.3580298485789791
$ws.Range("J6").Value = 0.3580298485789791
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7489546666666667
$ws.Range("N6").Value = 2.246864
$ws.Range("O6").Value = 0.05220789806691288
$ws.Range("P6").Value = 0.05220789806691287
$ws.Range("Q6").Value = 15.76787241614222
$ws.Range("R6").Value = 141.91085174528
$ws.Range("S6").Value = 0.01869198583952359
$ws.Range("T6").Value = 0.01869198583952359

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 21.05317333333333
$ws.Range("H7").Value = 63.15952
$ws.Range("I7").Value = 0.3580298485789791
$ws.Range("J7").Value = 0.3580298485789791
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.843693666666667
$ws.Range("N7").Value = 8.531081
$ws.Range("O7").Value = 0.1982273102638064
$ws.Range("P7").Value = 0.1982273102638064
$ws.Range("Q7").Value = 59.86877567123556
$ws.Range("R7").Value = 538.8189810411201
$ws.Range("S7").Value = 0.07097129387796894
$ws.Range("T7").Value = 0.07097129387796892

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 21.05317333333333
$ws.Range("H8").Value = 63.15952
$ws.Range("I8").Value = 0.3580298485789791
$ws.Range("J8").Value = 0.3580298485789791
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.337765
$ws.Range("N8").Value = 31.013295
$ws.Range("O8").Value = 0.7206216949842531
$ws.Range("P8").Value = 0.720621694984253
$ws.Range("Q8").Value = 217.6427584242666
$ws.Range("R8").Value = 1958.7848258184
$ws.Range("S8").Value = 0.2580040763379394
$ws.Range("T8").Value = 0.2580040763379394

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 21.05317333333333
$ws.Range("H9").Value = 63.15952
$ws.Range("I9").Value = 0.3580298485789791
$ws.Range("J9").Value = 0.3580298485789791
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4152066666666667
$ws.Range("N9").Value = 1.24562
$ws.Range("O9").Value = 0.02894309668502767
$ws.Range("P9").Value = 0.02894309668502767
$ws.Range("Q9").Value = 8.741417922488889
$ws.Range("R9").Value = 78.67276130239999
$ws.Range("S9").Value = 0.01036249252354721
$ws.Range("T9").Value = 0.01036249252354721

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 34.97741266666667
$ws.Range("H10").Value = 104.932238
$ws.Range("I10").Value = 0.5948251867999219
$ws.Range("J10").Value = 0.5948251867999219
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7489546666666667
$ws.Range("N10").Value = 2.246864
$ws.Range("O10").Value = 0.05220789806691288
$ws.Range("P10").Value = 0.05220789806691287
$ws.Range("Q10").Value = 26.19649644462578
$ws.Range("R10").Value = 235.768468001632
$ws.Range("S10").Value = 0.03105457272008274
$ws.Range("T10").Value = 0.03105457272008273

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 34.97741266666667
$ws.Range("H11").Value = 104.932238
$ws.Range("I11").Value = 0.5948251867999219
$ws.Range("J11").Value = 0.5948251867999219
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.843693666666667
$ws.Range("N11").Value = 8.531081
$ws.Range("O11").Value = 0.1982273102638064
$ws.Range("P11").Value = 0.1982273102638064
$ws.Range("Q11").Value = 99.46504687658647
$ws.Range("R11").Value = 895.1854218892781
$ws.Range("S11").Value = 0.1179105968565147
$ws.Range("T11").Value = 0.1179105968565147

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 34.97741266666667
$ws.Range("H12").Value = 104.932238
$ws.Range("I12").Value = 0.5948251867999219
$ws.Range("J12").Value = 0.5948251867999219
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 10.337765
$ws.Range("N12").Value = 31.013295
$ws.Range("O12").Value = 0.7206216949842531
$ws.Range("P12").Value = 0.720621694984253
$ws.Range("Q12").Value = 361.5882724560234
$ws.Range("R12").Value = 3254.29445210421
$ws.Range("S12").Value = 0.4286439343310847
$ws.Range("T12").Value = 0.4286439343310847

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 34.97741266666667
$ws.Range("H13").Value = 104.932238
$ws.Range("I13").Value = 0.5948251867999219
$ws.Range("J13").Value = 0.5948251867999219
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4152066666666667
$ws.Range("N13").Value = 1.24562
$ws.Range("O13").Value = 0.02894309668502767
$ws.Range("P13").Value = 0.02894309668502767
$ws.Range("Q13").Value = 14.52285492195111
$ws.Range("R13").Value = 130.70569429756
$ws.Range("S13").Value = 0.01721608289223979
$ws.Range("T13").Value = 0.01721608289223978

# Add new rows 14-17 for Resolving-Mac as sending cluster
# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Vtn"
$ws.Range("C14").Value = "Itga8"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.004706
$ws.Range("H14").Value = 0.014118
$ws.Range("I14").Value = 0.00008003014275976175315348293271000557069783098995685577392578125
$ws.Range("J14").Value = 0.00008003014275976175315348293271000557069783098995685577392578125
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.7489546666666667
$ws.Range("N14").Value = 2.246864
$ws.Range("O14").Value = 0.05220789806691288
$ws.Range("P14").Value = 0.05220789806691287
$ws.Range("Q14").Value = 0.003524580661333333
$ws.Range("R14").Value = 0.031721225952
$ws.Range("S14").Value = 0.0000041782055354821272751507377218871397417387925088405609130859375
$ws.Range("T14").Value = 0.000004178205535482125581084843213286461605093791149556636810302734375

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Vtn"
$ws.Range("C15").Value = "Itga8"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.004706
$ws.Range("H15").Value = 0.014118
$ws.Range("I15").Value = 0.00008003014275976175315348293271000557069783098995685577392578125
$ws.Range("J15").Value = 0.00008003014275976175315348293271000557069783098995685577392578125
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.843693666666667
$ws.Range("N15").Value = 8.531081
$ws.Range("O15").Value = 0.1982273102638064
$ws.Range("P15").Value = 0.1982273102638064
$ws.Range("Q15").Value = 0.01338242239533333
$ws.Range("R15").Value = 0.120441801558
$ws.Range("S15").Value = 0.00001586415993929601866920149622064428740486619062721729278564453125
$ws.Range("T15").Value = 0.0000158641599392960085048061291690402185849961824715137481689453125

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Vtn"
$ws.Range("C16").Value = "Itga8"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.004706
$ws.Range("H16").Value = 0.014118
$ws.Range("I16").Value = 0.00008003014275976175315348293271000557069783098995685577392578125
$ws.Range("J16").Value = 0.00008003014275976175315348293271000557069783098995685577392578125
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 10.337765
$ws.Range("N16").Value = 31.013295
$ws.Range("O16").Value = 0.7206216949842531
$ws.Range("P16").Value = 0.720621694984253
$ws.Range("Q16").Value = 0.04864952208999999
$ws.Range("R16").Value = 0.43784569881
$ws.Range("S16").Value = 0.000057671457125371258160655141278283508654567413032054901123046875
$ws.Range("T16").Value = 0.000057671457125371258160655141278283508654567413032054901123046875

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Vtn"
$ws.Range("C17").Value = "Itga8"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.004706
$ws.Range("H17").Value = 0.014118
$ws.Range("I17").Value = 0.00008003014275976175315348293271000557069783098995685577392578125
$ws.Range("J17").Value = 0.00008003014275976175315348293271000557069783098995685577392578125
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4152066666666667
$ws.Range("N17").Value = 1.24562
$ws.Range("O17").Value = 0.02894309668502767
$ws.Range("P17").Value = 0.02894309668502767
$ws.Range("Q17").Value = 0.001953962573333334
$ws.Range("R17").Value = 0.01758566316
$ws.Range("S17").Value = 0.00000231632015961235201309087287924182163578734616748988628387451171875
$ws.Range("T17").Value = 0.00000231632015961235201309087287924182163578734616748988628387451171875

